# Add new columns I (I0) and J (IF) to the sheet, matching the H column's
# header style, and populate the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold/centered/bordered) from H1 onto the
# two new header cells before writing their text, so they pick up the same
# cellXf ("s=1") as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(2, 8, 8),
    @(3, 8, 8),
    @(4, 8, 8),
    @(5, 6, 6),
    @(6, 8, 8),
    @(7, 6, 7),
    @(8, 8, 9),
    @(9, 6, 6),
    @(10, 6, 8),
    @(11, 7, 8),
    @(12, 9, 9),
    @(13, 8, 9),
    @(14, 5, 6),
    @(15, 8, 9),
    @(16, 7, 8),
    @(17, 9, 9),
    @(18, 8, 10),
    @(19, 7, 8),
    @(20, 7, 8),
    @(21, 5, 6),
    @(22, 6, 9),
    @(23, 5, 6),
    @(24, 6, 6),
    @(25, 1, 3),
    @(26, 1, 3),
    @(27, 3, 4)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 9).Value = $entry[1]
    $ws.Cells.Item($r, 10).Value = $entry[2]
}
